# "ppt project name change"
# Renames the placeholder project name to "Essay Grader" throughout the deck:
#   - Slide 1 title textbox: "Essay Meter" -> "Essay Grader"
#   - Slide 3 body: three "<project_name>" mentions collapse into a bold+underlined
#     "Essay Grader" run followed by the remaining sentence text.
#   - Slide 5 title: "Architecture of <project_name>" -> "Architecture of Essay Grader"

$p = $ppt.ActivePresentation

# --- Slide 1: "Essay Meter" -> "Essay Grader" -------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(2)
$titleShape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "Essay Grader"

# --- Slide 3: collapse each "<project_name>" mention into "Essay Grader" ----
$slide3 = $p.Slides.Item(3)
$body = $slide3.Shapes.Item(2)
$bodyTr = $body.TextFrame.TextRange

function Set-EssayGraderRun($para, $remainderText) {
    # Runs before: [1]="...<"  [2]="project_name"  [3]="> remainder"
    # Runs after:  [1]="Essay Grader" (bold, underline)  [2]="remainderText"
    $para.Runs(3,1).Text = $remainderText
    $para.Runs(2,1).Text = ""
    $para.Runs(1,1).Text = ""
    $anchor = $para.Runs(1,1)
    $newRun = $anchor.InsertBefore("Essay Grader")
    $newRun.Font.Bold = -1
    $newRun.Font.Underline = -1
}

Set-EssayGraderRun $bodyTr.Paragraphs(1,1) " is a product which gives linguistic support to all our users and enhance their English writing skill."
Set-EssayGraderRun $bodyTr.Paragraphs(3,1) " also provides an attractive dashboard to our users to keep track of their scores and their essays."
Set-EssayGraderRun $bodyTr.Paragraphs(11,1) " will have plagiarism feature to check the authenticity of the essay."

# --- Slide 5: "Architecture of <project_name>" -> "Architecture of Essay Grader"
$slide5 = $p.Slides.Item(5)
$titleShape5 = $slide5.Shapes.Item(1)
$para5 = $titleShape5.TextFrame.TextRange.Paragraphs(1,1)
$para5.Runs(3,1).Text = ""
$para5.Runs(2,1).Text = ""
$para5.Runs(1,1).Text = "Architecture of Essay Grader"
